# Insert 3 new rows before row 254 (shifts existing rows 254-344 down to 257-347)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("254:256").EntireRow.Insert()

# Populate the 3 new rows with the new weekly price entry
# (Terminal La Palmera de La Serena - Chirimoya, Cultivar IV Region, 2023-10-30,
#  $/bandeja 10 kilos, Provincia de Limari)

$newRowsData = @(
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", [datetime]"2023-10-30", 4, "Fruta", 100107, "Otros", 100107002, "Chirimoya", "Cultivar IV Región", "Especial", 500, 18000, 19000, 18500, "`$/bandeja 10 kilos", "Provincia de Limarí", 1850, 10),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", [datetime]"2023-10-30", 4, "Fruta", 100107, "Otros", 100107002, "Chirimoya", "Cultivar IV Región", "Primera", 600, 15000, 16000, 15500, "`$/bandeja 10 kilos", "Provincia de Limarí", 1550, 10),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", [datetime]"2023-10-30", 4, "Fruta", 100107, "Otros", 100107002, "Chirimoya", "Cultivar IV Región", "Segunda", 400, 12000, 13000, 12500, "`$/bandeja 10 kilos", "Provincia de Limarí", 1250, 10)
)

for ($i = 0; $i -lt 3; $i++) {
    $rowNum = 254 + $i
    $rowData = $newRowsData[$i]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowData[$col - 1]
    }
}
